$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 4) describing the "get courses detail by id" endpoint
$ws.Range("A4").Value = "Course"
$ws.Range("B4").Value = "get"
$ws.Range("C4").Value = "/courses"
$ws.Range("E4").Value = "id"
$ws.Range("G4").Value = "get courses detail by id"

# Update the selected cell to match the author's final cursor position
$ws.Range("G7").Select()
